$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''62.911.04'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.86%  '
$ws.Range("D3").Value = '''3.028.98'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.74%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''592.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("D6").Value = '''153.47'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.89%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '''3.023.90'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.70%  '
$ws.Range("D9").Value = '''0.514'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("D10").Value = '''6.90'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +15.46%  '
$ws.Range("E11").Value = '  +4.12%  '
$ws.Range("D12").Value = '''0.462'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.73%  '
$ws.Range("E13").Value = '  +3.28%  '
$ws.Range("D14").Value = '''35.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.61%  '
$ws.Range("E15").Value = '  -0.09%  '
$ws.Range("D16").Value = '''3.531.10'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.76%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").Value = '''7.08'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.49%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '''62.894.82'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.85%  '
$ws.Range("D19").Value = '''3.031.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.79%  '
$ws.Range("D20").Value = '''452.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.21%  '
$ws.Range("D21").Value = '''14.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.54%  '
$ws.Range("D22").Value = '''0.697'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.81%  '
$ws.Range("D23").Value = '''7.50'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.33%  '
$ws.Range("B24").Value = 'RenderToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D24").Value = '''11.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +11.52%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '''83.07'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.16%  '
$ws.Range("D26").Value = '''2.33'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.79%  '
$ws.Range("D27").Value = '''12.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.38%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").Value = '''7.52'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.51%  '
$ws.Range("E30").Value = '  +11.36%  '
$ws.Range("D31").Value = '''2.69'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.22%  '
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("D33").Value = '''27.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.95%  '
$ws.Range("E34").Value = '  +2.35%  '
$ws.Range("D35").Value = '''0.0₃0858'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.72%  '
$ws.Range("E36").Value = '  +2.50%  '
$ws.Range("D37").Value = '''5.91'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.21%  '
$ws.Range("D38").Value = '''3.14'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.18%  '
$ws.Range("E39").Value = '  +9.15%  '
$ws.Range("D40").Value = '''2.09'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.43%  '
$ws.Range("D41").Value = '''50.41'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.43%  '
$ws.Range("D42").Value = '''9.11'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.65%  '
$ws.Range("E43").Value = '  +16.58%  '
$ws.Range("D44").Value = '''44.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +14.61%  '
$ws.Range("D45").Value = '''392.69'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.34%  '
$ws.Range("D46").Value = '''0.0359'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.34%  '
$ws.Range("D47").Value = '''2.720.36'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.39%  '
$ws.Range("D48").Value = '''133.44'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.42%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '''25.45'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.83%  '
$ws.Range("B50").Value = 'USDe'
$ws.Range("C50").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D50").Value = '''1.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("E51").Value = '  +7.97%  '
